$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 11, leaving row 2 (2000年) and row 12 (2010年, which
# will shift up to row 3 after the deletion) so that the data year (2010) is
# the only remaining data row after row 1 (headers).
$ws.Range("A3:H11").EntireRow.Delete() | Out-Null

# Now the sheet has: row1=headers, row2=2000年(old), row3=2010年(old row12 data)
# Delete the old 2000年 row (now row 2), leaving only the 2010年 data in row 2.
$ws.Range("A2:H2").EntireRow.Delete() | Out-Null
